$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stash a copy of the bold/border/center header style (style index 1,
#    currently on A3) onto a scratch cell far outside the used range so it
#    survives the full-range clear below and can be re-applied afterwards
#    without Excel fabricating a brand-new (duplicate) style entry.
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy() | Out-Null
$ws.Range("AZ100").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2) Wipe the entire previously-used range (A1:AG19). We rebuild it from
#    scratch so the shared-string table comes out in exactly the order the
#    new workbook expects (first-write order == final <sst> order).
# ---------------------------------------------------------------------------
$ws.Range("A1:AG19").Clear()

# ---------------------------------------------------------------------------
# 3) Write all the new/kept label strings FIRST, in the precise order that
#    should end up in xl/sharedStrings.xml: column B top-to-bottom
#    (rows 2-29) and then row 2's remaining columns C-W left-to-right.
# ---------------------------------------------------------------------------
$colBValues = @(
  "HKL",
  "Spiral5",
  "RotRing OmegaMax-90",
  "Equal Angle",
  "Tilt Rotate",
  "CLR",
  "Rizzie Hex",
  "Thomas Hex",
  "Tilt Rotate_Partial",
  "RotRing OmegaMax-60",
  "Equal Angle_Partial",
  "Rizzie Hex_Partial",
  "ND Single",
  "RD Single",
  "TD Single",
  "Morris Single",
  "Ring Perpendicular to ND",
  "Ring Perpendicular to RD",
  "Ring Perpendicular to TD",
  "OffsetFTD",
  "OffsetATD",
  "OffsetF45",
  "OffsetA45",
  "OffsetFRD",
  "OffsetARD",
  "Gaussian Quadrature",
  "Michael-CCHex",
  "Michael-SNHex"
)

for ($i = 0; $i -lt $colBValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colBValues[$i]
}

$row2CWValues = @(
  "[4, 0, 0]",
  "[4, 2, 0]",
  "[3, 3, 3]",
  "[2, 2, 0]",
  "[2, 0, 0]",
  "[3, 1, 1]",
  "[3, 3, 1]",
  "[2, 2, 2]",
  "[1, 1, 1]",
  "[5, 1, 1]",
  "[4, 2, 2]",
  "1Pair-A",
  "1Pair-B",
  "2Pairs-A",
  "2Pairs-B",
  "3Pairs-A",
  "3Pairs-B",
  "3Pairs-C",
  "4Pairs",
  "5A4F",
  "MaxUnique"
)

for ($i = 0; $i -lt $row2CWValues.Length; $i++) {
    $ws.Cells.Item(2, $i + 3).Value = $row2CWValues[$i]
}

# ---------------------------------------------------------------------------
# 4) Fill in the numeric grid.
# ---------------------------------------------------------------------------

# Row 1: B1:W1 = 0..21
for ($c = 2; $c -le 23; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 2
}

# Column A: A2:A29 = 0..27
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Data body: C2:W29 (row2 already holds the bracket/pair header strings from
# step 3, so the numeric fill only applies to rows 3-29 here) = 1
for ($r = 3; $r -le 29; $r++) {
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# ---------------------------------------------------------------------------
# 5) Re-apply the stashed header style to the header row and header column.
# ---------------------------------------------------------------------------
$ws.Range("AZ100").Copy() | Out-Null
$ws.Range("B1:W1").PasteSpecial(-4122) | Out-Null
$ws.Range("AZ100").Copy() | Out-Null
$ws.Range("A2:A29").PasteSpecial(-4122) | Out-Null
$ws.Range("AZ100").Clear()

$excel.CutCopyMode = $false
